$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add new row of data: TF073410 / Core MDTF II
$ws.Range("A21").Value = "TF073410"
$ws.Range("B21").Value = "Core MDTF II"

# Copy style from row above (row 20) to new row 21 so formatting (shading) matches
$ws.Range("A20:B20").Copy()
$ws.Range("A21:B21").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Update the selection to match the new active cell recorded in the file
$ws.Range("B24").Select()
